# Apply the "Deploying to gh-pages" content update:
#  - Version bump 5.0.0 -> 6.0.0
#  - Date bump to the new publication date
#  - Publisher value filled in ("Alvearie Team")
#  - Duplicate "Contact / No display for ContactDetail" row replaced by a
#    single "Jurisdiction / United States of America" row
#  - Elements sheet "Short"/"Definition" for the root Extension row updated
#    to the resource's own Title/Description

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item(1)      # "Metadata" sheet
$elements = $wb.Worksheets.Item(2)  # "Elements" sheet

# The source data had an accidental duplicate of the "Contact" /
# "No display for ContactDetail" row (rows 10 and 11). Remove the extra
# row so the table below shifts up by one.
$meta.Rows.Item(10).Delete()

# Version
$meta.Range("B3").Value = "6.0.0"

# Date
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value (was blank)
$meta.Range("B9").Value = "Alvearie Team"

# The remaining "Contact" row becomes "Jurisdiction"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Elements sheet: root Extension row's Short/Definition now mirror the
# StructureDefinition's own Title/Description instead of the generic
# "Extension" / "An Extension" placeholders.
$elements.Range("K2").Value = "Local Number Of Units Per Service"
$elements.Range("L2").Value = "Customer-specific quantity of either services or units"
